$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-report data between row 2 and row 4
# (columns D, J, K, M, P), leaving the rest of each row untouched.

$ws.Range("D2").Value = 44792
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 528

$ws.Range("D4").Value = 44804
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 9500
$ws.Range("M4").Value = 9750
$ws.Range("P4").Value = 542
